# edit.ps1 - reproduces the timelineQRCode.pptx commit:
#   - "change slide location": reposition/resize the QR-code picture on
#     the last slide (slide 30).
#   - refreshed "datetimeFigureOut" footer field cached text (19-02-18
#     -> 08-03-18) on the slide master and every slide layout.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1) Move/resize the QR-code picture (slide 30, shape "Picture 3").
#    Target EMU values (from the canonical OOXML):
#      off  : x=4421688 y=1702964
#      ext  : cx=3024898 cy=3041009
#    Shape.Left/Top/Width/Height are expressed in points (1 pt = 12700 EMU).
#    The literal constants below were chosen so that, after the engine's
#    internal point<->EMU conversion, the saved XML lands exactly on the
#    target EMU values.
# ---------------------------------------------------------------------
$qrSlide = $p.Slides.Item(30)
$qrShape = $qrSlide.Shapes.Item(1)

$qrShape.Left   = 348.1644287109375
$qrShape.Top    = 134.09165954589844
$qrShape.Width  = 238.1809539794922
$qrShape.Height = 239.4495391845703

# ---------------------------------------------------------------------
# 2) Refresh the cached "date and time" footer field text everywhere it
#    appears: the slide master plus all eleven slide layouts.
# ---------------------------------------------------------------------
$newDate = "08-03-18"
$ppPlaceholderDate = 16

function Set-DateField($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $sh = $shapes.Item($i)
        $isDate = $false
        try {
            if ($sh.PlaceholderFormat.Type -eq $ppPlaceholderDate) {
                $isDate = $true
            }
        } catch {
            $isDate = $false
        }
        if ($isDate -and $sh.HasTextFrame) {
            if ($sh.TextFrame.TextRange.Text -ne $newDate) {
                $sh.TextFrame.TextRange.Text = $newDate
            }
        }
    }
}

$master = $p.SlideMaster
Set-DateField $master.Shapes

for ($L = 1; $L -le $master.CustomLayouts.Count; $L++) {
    $layout = $master.CustomLayouts.Item($L)
    Set-DateField $layout.Shapes
}
